$d = $word.ActiveDocument

$d.Content.Find.Execute("191÷4=47, 3", $true, $false, $false, $false, $false, $true, 1, $false, "837÷9=93, 0", 2) | Out-Null
$d.Content.Find.Execute("749÷8=93, 5", $true, $false, $false, $false, $false, $true, 1, $false, "768÷5=153, 3", 2) | Out-Null
$d.Content.Find.Execute("226÷4=56, 2", $true, $false, $false, $false, $false, $true, 1, $false, "140÷3=46, 2", 2) | Out-Null
$d.Content.Find.Execute("379÷3=126, 1", $true, $false, $false, $false, $false, $true, 1, $false, "878÷7=125, 3", 2) | Out-Null
$d.Content.Find.Execute("118÷2=59, 0", $true, $false, $false, $false, $false, $true, 1, $false, "292÷6=48, 4", 2) | Out-Null
$d.Content.Find.Execute("374÷5=74, 4", $true, $false, $false, $false, $false, $true, 1, $false, "634÷2=317, 0", 2) | Out-Null
$d.Content.Find.Execute("797÷7=113, 6", $true, $false, $false, $false, $false, $true, 1, $false, "172÷2=86, 0", 2) | Out-Null
$d.Content.Find.Execute("727÷6=121, 1", $true, $false, $false, $false, $false, $true, 1, $false, "767÷5=153, 2", 2) | Out-Null
$d.Content.Find.Execute("385÷7=55, 0", $true, $false, $false, $false, $false, $true, 1, $false, "347÷7=49, 4", 2) | Out-Null
$d.Content.Find.Execute("444÷8=55, 4", $true, $false, $false, $false, $false, $true, 1, $false, "833÷4=208, 1", 2) | Out-Null
$d.Content.Find.Execute("847÷9=94, 1", $true, $false, $false, $false, $false, $true, 1, $false, "383÷8=47, 7", 2) | Out-Null
$d.Content.Find.Execute("743÷3=247, 2", $true, $false, $false, $false, $false, $true, 1, $false, "326÷8=40, 6", 2) | Out-Null
$d.Content.Find.Execute("718÷2=359, 0", $true, $false, $false, $false, $false, $true, 1, $false, "519÷9=57, 6", 2) | Out-Null
$d.Content.Find.Execute("581÷9=64, 5", $true, $false, $false, $false, $false, $true, 1, $false, "429÷7=61, 2", 2) | Out-Null
$d.Content.Find.Execute("345÷6=57, 3", $true, $false, $false, $false, $false, $true, 1, $false, "971÷7=138, 5", 2) | Out-Null
$d.Content.Find.Execute("446÷9=49, 5", $true, $false, $false, $false, $false, $true, 1, $false, "540÷2=270, 0", 2) | Out-Null
$d.Content.Find.Execute("985÷5=197, 0", $true, $false, $false, $false, $false, $true, 1, $false, "480÷2=240, 0", 2) | Out-Null
$d.Content.Find.Execute("400÷6=66, 4", $true, $false, $false, $false, $false, $true, 1, $false, "209÷2=104, 1", 2) | Out-Null
$d.Content.Find.Execute("427÷5=85, 2", $true, $false, $false, $false, $false, $true, 1, $false, "987÷6=164, 3", 2) | Out-Null
$d.Content.Find.Execute("987÷7=141, 0", $true, $false, $false, $false, $false, $true, 1, $false, "479÷6=79, 5", 2) | Out-Null
$d.Content.Find.Execute("519÷5=103, 4", $true, $false, $false, $false, $false, $true, 1, $false, "438÷3=146, 0", 2) | Out-Null
$d.Content.Find.Execute("223÷3=74, 1", $true, $false, $false, $false, $false, $true, 1, $false, "499÷9=55, 4", 2) | Out-Null
$d.Content.Find.Execute("755÷3=251, 2", $true, $false, $false, $false, $false, $true, 1, $false, "322÷9=35, 7", 2) | Out-Null
$d.Content.Find.Execute("439÷5=87, 4", $true, $false, $false, $false, $false, $true, 1, $false, "907÷6=151, 1", 2) | Out-Null
$d.Content.Find.Execute("688÷3=229, 1", $true, $false, $false, $false, $false, $true, 1, $false, "254÷3=84, 2", 2) | Out-Null
